$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.448.14'
$ws.Range('E2').Value = '  +5.23%  '
$ws.Range('D3').Value = '2.055.26'
$ws.Range('E3').Value = '  +3.59%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.86'
$ws.Range('E5').Value = '  +3.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.654'
$ws.Range('E6').Value = '  +2.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '65.73'
$ws.Range('E7').Value = '  +13.81%  '
$ws.Range('E9').Value = '  +6.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.45'
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0770'
$ws.Range('E11').Value = '  +4.90%  '
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.927'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.91'
$ws.Range('E14').Value = '  +2.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.48'
$ws.Range('E15').Value = '  +25.80%  '
$ws.Range('D16').Value = '2.354.67'
$ws.Range('E16').Value = '  +3.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.58'
$ws.Range('E17').Value = '  +5.31%  '
$ws.Range('D18').Value = '2.057.07'
$ws.Range('E18').Value = '  +3.69%  '
$ws.Range('D19').Value = '37.351.75'
$ws.Range('E19').Value = '  +5.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.61'
$ws.Range('E20').Value = '  +3.23%  '
$ws.Range('E21').Value = '  +3.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.50'
$ws.Range('E22').Value = '  +6.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.16'
$ws.Range('E23').Value = '  +3.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.69'
$ws.Range('E24').Value = '  +4.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.41'
$ws.Range('E26').Value = '  +5.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.02'
$ws.Range('E27').Value = '  +9.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.27'
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.01'
$ws.Range('E29').Value = '  +4.10%  '
$ws.Range('E30').Value = '  +27.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.27'
$ws.Range('E31').Value = '  +8.40%  '
$ws.Range('E32').Value = '  +2.67%  '
$ws.Range('E33').Value = '  +9.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0628'
$ws.Range('E35').Value = '  +5.86%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.47'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.83'
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('E39').Value = '  +15.29%  '
$ws.Range('E40').Value = '  +35.46%  '
$ws.Range('E41').Value = '  +17.01%  '
$ws.Range('E42').Value = '  +4.60%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.80'
$ws.Range('E43').Value = '  +10.43%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.24'
$ws.Range('E44').Value = '  +1.67%  '
$ws.Range('E45').Value = '  +6.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0219'
$ws.Range('E46').Value = '  +3.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '97.06'
$ws.Range('E47').Value = '  +5.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.93'
$ws.Range('E48').Value = '  +3.13%  '
$ws.Range('D49').Value = '1.419.01'
$ws.Range('E49').Value = '  +3.31%  '
$ws.Range('E50').Value = '  +1.72%  '
$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.76'
$ws.Range('E51').Value = '  +7.67%  '
